# Applies the "456a3b4" content refresh described in the commit:
#  - Inserts a new event ("广州·第五人格ONLY") into both the "展览" sheet
#    (row 26) and the "全部类型" sheet (row 34), pushing later rows down.
#  - Bumps several "想去人数" (F column) vote counts that changed between
#    scrapes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data: independent F-column (want-to-go count) bumps that are
# NOT related to the row insertion - same events, updated vote counts.
# Keyed by the row number *before* any insertion happens on that sheet.
# ---------------------------------------------------------------------

function Update-FCounts($ws, $map) {
    foreach ($row in $map.Keys) {
        $ws.Cells.Item($row, 6).Value = $map[$row]
    }
}

function Insert-Row26Event($ws, $insertRow) {
    $ws.Rows.Item($insertRow).Insert()

    # Borrow the index-column formatting (bold / centered / bordered)
    # from the row just above so the new A cell matches its neighbours.
    $ws.Cells.Item($insertRow - 1, 1).Copy()
    $ws.Cells.Item($insertRow, 1).PasteSpecial(-4122)

    $ws.Cells.Item($insertRow, 2).NumberFormat = "@"
    $ws.Cells.Item($insertRow, 2).Value = "2024-06-23"
    $ws.Cells.Item($insertRow, 3).Value = "广州·第五人格ONLY"
    $ws.Cells.Item($insertRow, 4).Value = "奥体南路12号 优托邦(奥体旗舰店)"
    $ws.Cells.Item($insertRow, 5).Value = "2024.06.23 10:00-06.23 17:00"
    $ws.Cells.Item($insertRow, 6).Value = 1
    $ws.Cells.Item($insertRow, 7).Value = 55
    $ws.Cells.Item($insertRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86276"
    $ws.Cells.Item($insertRow, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/c9J2CXNZ1716259888686.jpeg"
}

# ===================== Sheet 1: "展览" =====================
$ws1 = $wb.Worksheets.Item(1)

# Independent "want to go" count bumps (rows are pre-insertion numbers).
$f1 = @{
    2  = 26
    3  = 1802
    4  = 467
    8  = 349
    9  = 1756
    10 = 377
    14 = 692
    15 = 12890
    16 = 12862
    17 = 962
    20 = 524
    21 = 53
    22 = 588
}
Update-FCounts $ws1 $f1

# Insert the new event as row 26 (existing rows 26-29 become 27-30).
Insert-Row26Event $ws1 26

# Re-sequence the index column (A) so it keeps tracking "row - 1" for
# the rows that shifted down, and bump the two vote counts that also
# changed for the events now sitting in rows 28-30.
$ws1.Cells.Item(26, 1).Value = 25
$ws1.Cells.Item(27, 1).Value = 26
$ws1.Cells.Item(28, 1).Value = 27
$ws1.Cells.Item(28, 6).Value = 93
$ws1.Cells.Item(29, 1).Value = 28
$ws1.Cells.Item(29, 6).Value = 258
$ws1.Cells.Item(30, 1).Value = 29
$ws1.Cells.Item(30, 6).Value = 686

# ===================== Sheet 4: "全部类型" =====================
$ws4 = $wb.Worksheets.Item(4)

$f4 = @{
    4  = 26
    5  = 1802
    6  = 467
    13 = 349
    14 = 1756
    15 = 377
    20 = 692
    21 = 12890
    22 = 12862
    23 = 962
    26 = 524
    27 = 53
    28 = 588
}
Update-FCounts $ws4 $f4

# Insert the new event as row 34 (existing rows 34-43 become 35-44).
Insert-Row26Event $ws4 34

# Re-sequence the index column (A) for the shifted rows and bump the
# two vote counts that also changed for the events now in rows 38-40.
$ws4.Cells.Item(34, 1).Value = 33
$ws4.Cells.Item(35, 1).Value = 34
$ws4.Cells.Item(36, 1).Value = 35
$ws4.Cells.Item(37, 1).Value = 36
$ws4.Cells.Item(38, 1).Value = 37
$ws4.Cells.Item(38, 6).Value = 93
$ws4.Cells.Item(39, 1).Value = 38
$ws4.Cells.Item(39, 6).Value = 258
$ws4.Cells.Item(40, 1).Value = 39
$ws4.Cells.Item(40, 6).Value = 686
$ws4.Cells.Item(41, 1).Value = 40
$ws4.Cells.Item(42, 1).Value = 41
$ws4.Cells.Item(43, 1).Value = 42
$ws4.Cells.Item(44, 1).Value = 43

Write-Output "edit complete"
